$wb = $excel.ActiveWorkbook

# --- 1) "sets" sheet: E4 changes from 2 to 4 ---
$sets = $wb.Worksheets.Item("sets")
$sets.Range("E4").Value = 4

# --- 2) "rallies" sheet: append two new rows (90 and 91) ---
$rallies = $wb.Worksheets.Item("rallies")

# Row 90
$rallies.Cells.Item(90, 1).Value = 89
$rallies.Cells.Item(90, 2).Value = 1
$rallies.Cells.Item(90, 3).Value = 3
$rallies.Cells.Item(90, 4).Value = 27
$rallies.Cells.Item(90, 5).Value = "ADV"
$rallies.Cells.Item(90, 6).Value = "ADVERSÁRIO"
$rallies.Cells.Item(90, 7).Value = ""
$rallies.Cells.Item(90, 8).Value = "PIPE"
$rallies.Cells.Item(90, 9).Value = "PONTO"
$rallies.Cells.Item(90, 10).Value = "ADV"
$rallies.Cells.Item(90, 11).Value = 24
$rallies.Cells.Item(90, 12).Value = 3
$rallies.Cells.Item(90, 13).Value = "0  pi"
$rallies.Cells.Item(90, 14).Value = "FRENTE"
$rallies.Cells.Item(90, 15).Value = "FRENTE"
$rallies.Cells.Item(90, 16).Value = "FRENTE"

# Row 91
$rallies.Cells.Item(91, 1).Value = 90
$rallies.Cells.Item(91, 2).Value = 1
$rallies.Cells.Item(91, 3).Value = 3
$rallies.Cells.Item(91, 4).Value = 28
$rallies.Cells.Item(91, 5).Value = "ADV"
$rallies.Cells.Item(91, 6).Value = "ADVERSÁRIO"
$rallies.Cells.Item(91, 7).Value = ""
$rallies.Cells.Item(91, 8).Value = "RECEPÇÃO"
$rallies.Cells.Item(91, 9).Value = "PONTO"
$rallies.Cells.Item(91, 10).Value = "ADV"
$rallies.Cells.Item(91, 11).Value = 24
$rallies.Cells.Item(91, 12).Value = 4
$rallies.Cells.Item(91, 13).Value = "0  re"
$rallies.Cells.Item(91, 14).Value = "FRENTE"
$rallies.Cells.Item(91, 15).Value = "FRENTE"
$rallies.Cells.Item(91, 16).Value = "FRENTE"
